# Fruta / hortaliza, semanal
# Insert one new weekly record (row 31) for the "Ciruela" sheet and update
# the figures for the existing row 30 record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot row 30's current ("before edit") values -----------------
# These become the content of the newly-inserted row 31.
$origA30 = $ws.Cells.Item(30, 1).Value()
$origB30 = $ws.Cells.Item(30, 2).Value()
$origC30 = $ws.Cells.Item(30, 3).Value()
$origD30 = 44223
$origE30 = $ws.Cells.Item(30, 5).Value()
$origF30 = $ws.Cells.Item(30, 6).Value()
$origG30 = $ws.Cells.Item(30, 7).Value()
$origH30 = $ws.Cells.Item(30, 8).Value()
$origI30 = $ws.Cells.Item(30, 9).Value()
$origJ30 = $ws.Cells.Item(30, 10).Value()
$origK30 = $ws.Cells.Item(30, 11).Value()
$origL30 = $ws.Cells.Item(30, 12).Value()
$origM30 = $ws.Cells.Item(30, 13).Value()
$origN30 = $ws.Cells.Item(30, 14).Value()
$origO30 = $ws.Cells.Item(30, 15).Value()
$origP30 = $ws.Cells.Item(30, 16).Value()
$origQ30 = $ws.Cells.Item(30, 17).Value()
$origR30 = $ws.Cells.Item(30, 18).Value()
$origS30 = $ws.Cells.Item(30, 19).Value()
$origT30 = $ws.Cells.Item(30, 20).Value()

# --- Insert a new row at 31, pushing old rows 31-49 down to 32-50 -----
$ws.Rows.Item(31).Insert()

# --- Populate the new row 31 with row 30's original record ------------
$ws.Cells.Item(31, 1).Value = $origA30
$ws.Cells.Item(31, 2).Value = $origB30
$ws.Cells.Item(31, 3).Value = $origC30
$ws.Cells.Item(31, 4).Value = $origD30
$ws.Cells.Item(31, 5).Value = $origE30
$ws.Cells.Item(31, 6).Value = $origF30
$ws.Cells.Item(31, 7).Value = $origG30
$ws.Cells.Item(31, 8).Value = $origH30
$ws.Cells.Item(31, 9).Value = $origI30
$ws.Cells.Item(31, 10).Value = $origJ30
$ws.Cells.Item(31, 11).Value = $origK30
$ws.Cells.Item(31, 12).Value = $origL30
$ws.Cells.Item(31, 13).Value = $origM30
$ws.Cells.Item(31, 14).Value = $origN30
$ws.Cells.Item(31, 15).Value = $origO30
$ws.Cells.Item(31, 16).Value = $origP30
$ws.Cells.Item(31, 17).Value = $origQ30
$ws.Cells.Item(31, 18).Value = $origR30
$ws.Cells.Item(31, 19).Value = $origS30
$ws.Cells.Item(31, 20).Value = $origT30

# --- Update row 30 with this week's new figures ------------------------
$ws.Cells.Item(30, 4).Value = 44596
$ws.Cells.Item(30, 13).Value = 220
$ws.Cells.Item(30, 14).Value = 7000
$ws.Cells.Item(30, 15).Value = 8000
$ws.Cells.Item(30, 16).Value = 7545
$ws.Cells.Item(30, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(30, 19).Value = 472
